$wb = $excel.ActiveWorkbook

# Update Status on every sheet (shared string used across Overview/zh-cn/de-de)
foreach ($name in @("Overview", "zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            if ($cell.Value() -eq "Ready for handoff") {
                $cell.Value = "Handoff transform failed"
            }
        }
    }
}

foreach ($name in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2: clear the "Latest Handoff File" hyperlink/value, reset the handoff
    # datetime, and flip the handoff reason to "Ignored" (handoff failed).
    $ws.Range("C2").Hyperlinks.Delete()
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"
}
